# Applies the "new paper pulp input file" change:
#  - adds a "seasonal_efficiency" column (E) to the year_Vecteurs sheet
#  - boxes the whole table in thin borders, keeps the hydrogen rows highlighted
#  - updates the selections / active sheet left behind by the editing session

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# year_Vecteurs: new "seasonal_efficiency" column + formatting
# ---------------------------------------------------------------------------
$wsVec = $wb.Worksheets.Item("year_Vecteurs")
[void]$wsVec.Activate()

$wsVec.Range("E1").Value = "seasonal_efficiency"

$wsVec.Range("E2").Value = 0.5588
$wsVec.Range("E3").Value = 0.4536
$wsVec.Range("E4").Value = 1
$wsVec.Range("E5").Value = 1
$wsVec.Range("E6").Value = 1
$wsVec.Range("E7").Value = 0.4536
$wsVec.Range("E8").Value = 0.66612401833447832
$wsVec.Range("E9").Value = 0.53659990365832977
$wsVec.Range("E10").Value = 1
$wsVec.Range("E11").Value = 1
$wsVec.Range("E12").Value = 1
$wsVec.Range("E13").Value = 0.53659990365832977

# Box every data cell (A1:E13) in a thin border. The existing yellow
# "hydrogen" rows (C7:D7 / C13:D13) already carry their fill, so this alone
# gives them fill+border while the rest get border-only.
$wsVec.Range("A1:E13").Borders.LineStyle = 1

[void]$wsVec.Range("H9").Select()

# ---------------------------------------------------------------------------
# Production_system_year: selection left on A1:A5 (no single active cell)
# ---------------------------------------------------------------------------
$wsYear = $wb.Worksheets.Item("Production_system_year")
[void]$wsYear.Activate()
[void]$wsYear.Range("A1:A5").Select()

# ---------------------------------------------------------------------------
# retrofit_Transition: selection moved to D16
# ---------------------------------------------------------------------------
$wsRetrofit = $wb.Worksheets.Item("retrofit_Transition")
[void]$wsRetrofit.Activate()
[void]$wsRetrofit.Range("D16").Select()

# ---------------------------------------------------------------------------
# 0D: selection moved to A4, ends up the active / tab-selected sheet
# ---------------------------------------------------------------------------
$ws0D = $wb.Worksheets.Item("0D")
[void]$ws0D.Activate()
[void]$ws0D.Range("A4").Select()
